$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Version" column (F) header.
$ws.Range("F1").Value = "Version"

# Populate distinct version values in the order that first introduces
# each new shared string, matching the target workbook's string table.
$ws.Range("F4").Value = "3.4.1.113"
$ws.Range("F5").Value = "3.3.1.110"
$ws.Range("F6").Value = "3.5.2.114"
$ws.Range("F10").Value = "3.6.0.0"
$ws.Range("F2").Value = "alfa1"
$ws.Range("F3").Value = "alfa1"

# Remaining rows reuse already-registered shared strings.
$ws.Range("F7").Value = "3.5.2.114"
$ws.Range("F8").Value = "3.5.2.114"
$ws.Range("F9").Value = "3.5.2.114"
$ws.Range("F11").Value = "3.5.2.114"

# Leave the active selection where editing stopped (just below the data).
[void]$ws.Range("A12").Select()
